$wb = $excel.ActiveWorkbook

# Work on the "On-Site" sheet: remove the "Control" row (row 4)
$onSite = $wb.Worksheets.Item("On-Site")
[void]$onSite.Activate()
[void]$onSite.Rows.Item(4).Delete()

# Select the row that now occupies row 4 (previously row 5 - "Support Structure")
[void]$onSite.Range("A4:XFD4").Select()

# Remove the entire "Replacement" sheet/worksheet from the workbook
$excel.DisplayAlerts = $false
[void]$wb.Worksheets.Item("Replacement").Delete()
$excel.DisplayAlerts = $true

# Keep "On-Site" as the active sheet
[void]$onSite.Activate()
